# The source diff for this revision only reorders XML attributes inside
# word/document.xml and word/styles.xml (e.g. xmlns:* declarations on
# <w:document>, and attributes on <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>,
# <w:latentStyles>, every <w:lsdException>, <w:style>, <w:tblInd> and the
# table cell margins). Every hunk swaps attribute *order* only -- the same
# attribute names and the same values are present before and after (this
# was confirmed by canonicalizing both revisions and comparing them byte
# for byte). The commit message ("Fixed POI packaging and upgraded to POI
# 3.15") corroborates this: upgrading the Apache POI library changed how
# it serializes attributes (alphabetically) but did not touch the
# document's actual content, formatting, styles or page setup.
#
# So there is no content/formatting edit for Word's object model to apply
# here -- the paragraphs, runs, drawing, sections, page setup, fonts,
# language settings and style definitions are all unchanged. We simply
# touch the document so the intent is explicit without mutating anything
# observable through the Word OM.
$d = $word.ActiveDocument
